$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 25,14
$data[0,0] = 1.444287666666667
$data[0,1] = 4.332863
$data[0,2] = 0.006189216566550864
$data[0,3] = 0.006202528009901729
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 121.928739
$data[0,7] = 365.786217
$data[0,8] = 0.2282232151508951
$data[0,9] = 0.2419720431319445
$data[0,10] = 176.100173949919
$data[0,11] = 1584.901565549271
$data[0,12] = 0.001412522904083422
$data[0,13] = 0.001500838375139035
$data[1,0] = 1.444287666666667
$data[1,1] = 4.332863
$data[1,2] = 0.006189216566550864
$data[1,3] = 0.006202528009901729
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 147.91433
$data[1,7] = 443.74299
$data[1,8] = 0.2768624053389947
$data[1,9] = 0.2935413991166814
$data[1,10] = 213.6308425422633
$data[1,11] = 1922.67758288037
$data[1,12] = 0.001713561385779226
$data[1,13] = 0.001820698750086959
$data[2,0] = 1.444287666666667
$data[2,1] = 4.332863
$data[2,2] = 0.006189216566550864
$data[2,3] = 0.006202528009901729
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 83.50496933333334
$data[2,7] = 250.514908
$data[2,8] = 0.1563025480180701
$data[2,9] = 0.1657186665504434
$data[2,10] = 120.6051973135116
$data[2,11] = 1085.446775821604
$data[2,12] = 0.000967390319587551
$data[2,13] = 0.00102787467104269
$data[3,0] = 1.444287666666667
$data[3,1] = 4.332863
$data[3,2] = 0.006189216566550864
$data[3,3] = 0.006202528009901729
$data[3,4] = 2
$data[3,5] = 1
$data[3,6] = 91.06846250000001
$data[3,7] = 182.136925
$data[3,8] = 0.1704597085236707
$data[3,9] = 0.1204857969594293
$data[3,10] = 131.5290572110459
$data[3,11] = 789.1743432662751
$data[3,12] = 0.001055012051924134
$data[3,13] = 0.0007473165304361931
$data[4,0] = 1.444287666666667
$data[4,1] = 4.332863
$data[4,2] = 0.006189216566550864
$data[4,3] = 0.006202528009901729
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 89.83562999999999
$data[4,7] = 269.50689
$data[4,8] = 0.1681521229683693
$data[4,9] = 0.1782820942415013
$data[4,10] = 129.74849243623
$data[4,11] = 1167.73643192607
$data[4,12] = 0.001040729905176529
$data[4,13] = 0.001105799683196851
$data[5,0] = 0.7051769999999999
$data[5,1] = 2.115531
$data[5,2] = 0.003021900187532335
$data[5,3] = 0.003028399532437424
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 121.928739
$data[5,7] = 365.786217
$data[5,8] = 0.2282232151508951
$data[5,9] = 0.2419720431319445
$data[5,10] = 85.98134238180299
$data[5,11] = 773.8320814362269
$data[5,12] = 0.0006896677766637222
$data[5,13] = 0.0007327880222837091
$data[6,0] = 0.7051769999999999
$data[6,1] = 2.115531
$data[6,2] = 0.003021900187532335
$data[6,3] = 0.003028399532437424
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 147.91433
$data[6,7] = 443.74299
$data[6,8] = 0.2768624053389947
$data[6,9] = 0.2935413991166814
$data[6,10] = 104.30578348641
$data[6,11] = 938.7520513776899
$data[6,12] = 0.0008366505546145614
$data[6,13] = 0.0008889606358359854
$data[7,0] = 0.7051769999999999
$data[7,1] = 2.115531
$data[7,2] = 0.003021900187532335
$data[7,3] = 0.003028399532437424
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 83.50496933333334
$data[7,7] = 250.514908
$data[7,8] = 0.1563025480180701
$data[7,9] = 0.1657186665504434
$data[7,10] = 58.885783759572
$data[7,11] = 529.9720538361479
$data[7,12] = 0.0004723306991675877
$data[7,13] = 0.0005018623322975162
$data[8,0] = 0.7051769999999999
$data[8,1] = 2.115531
$data[8,2] = 0.003021900187532335
$data[8,3] = 0.003028399532437424
$data[8,4] = 2
$data[8,5] = 1
$data[8,6] = 91.06846250000001
$data[8,7] = 182.136925
$data[8,8] = 0.1704597085236707
$data[8,9] = 0.1204857969594293
$data[8,10] = 64.21938518036251
$data[8,11] = 385.316311082175
$data[8,12] = 0.0005151122251543877
$data[8,13] = 0.0003648791311772862
$data[9,0] = 0.7051769999999999
$data[9,1] = 2.115531
$data[9,2] = 0.003021900187532335
$data[9,3] = 0.003028399532437424
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 89.83562999999999
$data[9,7] = 269.50689
$data[9,8] = 0.1681521229683693
$data[9,9] = 0.1782820942415013
$data[9,10] = 63.35002005650999
$data[9,11] = 570.15018050859
$data[9,12] = 0.0005081389319320753
$data[9,13] = 0.0005399094108429273
$data[10,0] = 123.254125
$data[10,1] = 369.762375
$data[10,2] = 0.528181808895687
$data[10,3] = 0.5293177947110922
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 121.928739
$data[10,7] = 365.786217
$data[10,8] = 0.2282232151508951
$data[10,9] = 0.2419720431319445
$data[10,10] = 15028.22003779837
$data[10,11] = 135253.9803401854
$data[10,12] = 0.1205433506103893
$data[10,13] = 0.1280801082523382
$data[11,0] = 123.254125
$data[11,1] = 369.762375
$data[11,2] = 0.528181808895687
$data[11,3] = 0.5293177947110922
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 147.91433
$data[11,7] = 443.74299
$data[11,8] = 0.2768624053389947
$data[11,9] = 0.2935413991166814
$data[11,10] = 18231.05131911125
$data[11,11] = 164079.4618720012
$data[11,12] = 0.1462336860671611
$data[11,13] = 0.1553766860368504
$data[12,0] = 123.254125
$data[12,1] = 369.762375
$data[12,2] = 0.528181808895687
$data[12,3] = 0.5293177947110922
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 83.50496933333334
$data[12,7] = 250.514908
$data[12,8] = 0.1563025480180701
$data[12,9] = 0.1657186665504434
$data[12,10] = 10292.33192833183
$data[12,11] = 92630.9873549865
$data[12,12] = 0.08255616254718923
$data[12,13] = 0.08771783912094352
$data[13,0] = 123.254125
$data[13,1] = 369.762375
$data[13,2] = 0.528181808895687
$data[13,3] = 0.5293177947110922
$data[13,4] = 2
$data[13,5] = 1
$data[13,6] = 91.06846250000001
$data[13,7] = 182.136925
$data[13,8] = 0.1704597085236707
$data[13,9] = 0.1204857969594293
$data[13,10] = 11224.56366053281
$data[13,11] = 67347.38196319688
$data[13,12] = 0.09003371719186397
$data[13,13] = 0.06377527634057355
$data[14,0] = 123.254125
$data[14,1] = 369.762375
$data[14,2] = 0.528181808895687
$data[14,3] = 0.5293177947110922
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 89.83562999999999
$data[14,7] = 269.50689
$data[14,8] = 0.1681521229683693
$data[14,9] = 0.1782820942415013
$data[14,10] = 11072.61196947375
$data[14,11] = 99653.50772526376
$data[14,12] = 0.08881489247908329
$data[14,13] = 0.09436788496038657
$data[15,0] = 1.502435
$data[15,1] = 3.00487
$data[15,2] = 0.006438395761993292
$data[15,3] = 0.00430149541795192
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 121.928739
$data[15,7] = 365.786217
$data[15,8] = 0.2282232151508951
$data[15,9] = 0.2419720431319445
$data[15,10] = 183.190004979465
$data[15,11] = 1099.14002987679
$data[15,12] = 0.001469391381216006
$data[15,13] = 0.001040841634804524
$data[16,0] = 1.502435
$data[16,1] = 3.00487
$data[16,2] = 0.006438395761993292
$data[16,3] = 0.00430149541795192
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 147.91433
$data[16,7] = 443.74299
$data[16,8] = 0.2768624053389947
$data[16,9] = 0.2935413991166814
$data[16,10] = 222.23166639355
$data[16,11] = 1333.3899983613
$data[16,12] = 0.001782549737189852
$data[16,13] = 0.001262666983279601
$data[17,0] = 1.502435
$data[17,1] = 3.00487
$data[17,2] = 0.006438395761993292
$data[17,3] = 0.00430149541795192
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 83.50496933333334
$data[17,7] = 250.514908
$data[17,8] = 0.1563025480180701
$data[17,9] = 0.1657186665504434
$data[17,10] = 125.4607886003267
$data[17,11] = 752.76473160196
$data[17,12] = 0.001006337662748295
$data[17,13] = 0.0007128380848358344
$data[18,0] = 1.502435
$data[18,1] = 3.00487
$data[18,2] = 0.006438395761993292
$data[18,3] = 0.00430149541795192
$data[18,4] = 2
$data[18,5] = 1
$data[18,6] = 91.06846250000001
$data[18,7] = 182.136925
$data[18,8] = 0.1704597085236707
$data[18,9] = 0.1204857969594293
$data[18,10] = 136.8244454561875
$data[18,11] = 547.29778182475
$data[18,12] = 0.001097487064949414
$data[18,13] = 0.0005182691035492706
$data[19,0] = 1.502435
$data[19,1] = 3.00487
$data[19,2] = 0.006438395761993292
$data[19,3] = 0.00430149541795192
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 89.83562999999999
$data[19,7] = 269.50689
$data[19,8] = 0.1681521229683693
$data[19,9] = 0.1782820942415013
$data[19,10] = 134.97219475905
$data[19,11] = 809.8331685543
$data[19,12] = 0.001082629915889724
$data[19,13] = 0.0007668796114826903
$data[20,0] = 106.4494656666667
$data[20,1] = 319.348397
$data[20,2] = 0.4561686785882365
$data[20,3] = 0.4571497823286167
$data[20,4] = 3
$data[20,5] = 1
$data[20,6] = 121.928739
$data[20,7] = 365.786217
$data[20,8] = 0.2282232151508951
$data[20,9] = 0.2419720431319445
$data[20,10] = 12979.24911596046
$data[20,11] = 116813.2420436441
$data[20,12] = 0.1041082824785426
$data[20,13] = 0.1106174668473791
$data[21,0] = 106.4494656666667
$data[21,1] = 319.348397
$data[21,2] = 0.4561686785882365
$data[21,3] = 0.4571497823286167
$data[21,4] = 3
$data[21,5] = 1
$data[21,6] = 147.91433
$data[21,7] = 443.74299
$data[21,8] = 0.2768624053389947
$data[21,9] = 0.2935413991166814
$data[21,10] = 15745.401392943
$data[21,11] = 141708.612536487
$data[21,12] = 0.1262959575942499
$data[21,13] = 0.1341923867106285
$data[22,0] = 106.4494656666667
$data[22,1] = 319.348397
$data[22,2] = 0.4561686785882365
$data[22,3] = 0.4571497823286167
$data[22,4] = 3
$data[22,5] = 1
$data[22,6] = 83.50496933333334
$data[22,7] = 250.514908
$data[22,8] = 0.1563025480180701
$data[22,9] = 0.1657186665504434
$data[22,10] = 8889.059366044719
$data[22,11] = 80001.53429440247
$data[22,12] = 0.07130032678937739
$data[22,13] = 0.07575825234132381
$data[23,0] = 106.4494656666667
$data[23,1] = 319.348397
$data[23,2] = 0.4561686785882365
$data[23,3] = 0.4571497823286167
$data[23,4] = 2
$data[23,5] = 1
$data[23,6] = 91.06846250000001
$data[23,7] = 182.136925
$data[23,8] = 0.1704597085236707
$data[23,9] = 0.1204857969594293
$data[23,10] = 9694.18917220987
$data[23,11] = 58165.13503325923
$data[23,12] = 0.07775837998977883
$data[23,13] = 0.05508005585369304
$data[24,0] = 106.4494656666667
$data[24,1] = 319.348397
$data[24,2] = 0.4561686785882365
$data[24,3] = 0.4571497823286167
$data[24,4] = 3
$data[24,5] = 1
$data[24,6] = 89.83562999999999
$data[24,7] = 269.50689
$data[24,8] = 0.1681521229683693
$data[24,9] = 0.1782820942415013
$data[24,10] = 9562.954811328369
$data[24,11] = 86066.59330195532
$data[24,12] = 0.07670573173628767
$data[24,13] = 0.08150162057559225

$ws.Range("G2:T26").Value = $data
Write-Output "Done setting G2:T26"